$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy the date-format style from A52 into the newly added A53 cell
$ws.Range("A52").Copy($ws.Range("A53"))

# Rewrite all data rows (2-53) with the corrected / newly simulated values
$ws.Cells.Item(2, 1).Value = 39400
$ws.Cells.Item(2, 2).Value = 2007
$ws.Cells.Item(2, 3).Value = 2.070003986395053
$ws.Cells.Item(2, 4).Value = 2008
$ws.Cells.Item(2, 5).Value = 0.6967455006573253

$ws.Cells.Item(3, 1).Value = 39583
$ws.Cells.Item(3, 2).Value = 2008
$ws.Cells.Item(3, 3).Value = 0.9380533998416762
$ws.Cells.Item(3, 4).Value = 2009
$ws.Cells.Item(3, 5).Value = 1.155022846258058

$ws.Cells.Item(4, 1).Value = 39765
$ws.Cells.Item(4, 2).Value = 2008
$ws.Cells.Item(4, 3).Value = 0.517569958955022
$ws.Cells.Item(4, 4).Value = 2009
$ws.Cells.Item(4, 5).Value = -0.6367039903686034

$ws.Cells.Item(5, 1).Value = 39948
$ws.Cells.Item(5, 2).Value = 2009
$ws.Cells.Item(5, 3).Value = -4.857241224140929
$ws.Cells.Item(5, 4).Value = 2010
$ws.Cells.Item(5, 5).Value = -1.691674259276643

$ws.Cells.Item(6, 1).Value = 40130
$ws.Cells.Item(6, 2).Value = 2009
$ws.Cells.Item(6, 3).Value = -3.956152295564896
$ws.Cells.Item(6, 4).Value = 2010
$ws.Cells.Item(6, 5).Value = -0.6950853968889392

$ws.Cells.Item(7, 1).Value = 40310
$ws.Cells.Item(7, 2).Value = 2010
$ws.Cells.Item(7, 3).Value = 0.3625742673738941
$ws.Cells.Item(7, 4).Value = 2011
$ws.Cells.Item(7, 5).Value = -1.432000573345915

$ws.Cells.Item(8, 1).Value = 40494
$ws.Cells.Item(8, 2).Value = 2010
$ws.Cells.Item(8, 3).Value = 1.234995474941392
$ws.Cells.Item(8, 4).Value = 2011
$ws.Cells.Item(8, 5).Value = -0.2098161877568061

$ws.Cells.Item(9, 1).Value = 40676
$ws.Cells.Item(9, 2).Value = 2011
$ws.Cells.Item(9, 3).Value = 0.9965309787904442
$ws.Cells.Item(9, 4).Value = 2012
$ws.Cells.Item(9, 5).Value = -0.02552935725171901

$ws.Cells.Item(10, 1).Value = 40862
$ws.Cells.Item(10, 2).Value = 2011
$ws.Cells.Item(10, 3).Value = 0.899360810820804
$ws.Cells.Item(10, 4).Value = 2012
$ws.Cells.Item(10, 5).Value = 1.205741443109987

$ws.Cells.Item(11, 1).Value = 41044
$ws.Cells.Item(11, 2).Value = 2012
$ws.Cells.Item(11, 3).Value = 0.6836026627130565
$ws.Cells.Item(11, 4).Value = 2013
$ws.Cells.Item(11, 5).Value = 1.279644059586871

$ws.Cells.Item(12, 1).Value = 41228
$ws.Cells.Item(12, 2).Value = 2012
$ws.Cells.Item(12, 3).Value = 0.9010266119894084
$ws.Cells.Item(12, 4).Value = 2013
$ws.Cells.Item(12, 5).Value = 1.531699207045123

$ws.Cells.Item(13, 1).Value = 41409
$ws.Cells.Item(13, 2).Value = 2013
$ws.Cells.Item(13, 3).Value = 0.2336505480021955
$ws.Cells.Item(13, 4).Value = 2014
$ws.Cells.Item(13, 5).Value = 0.4758549421990166

$ws.Cells.Item(14, 1).Value = 41592
$ws.Cells.Item(14, 2).Value = 2013
$ws.Cells.Item(14, 3).Value = 0.02019328874804938
$ws.Cells.Item(14, 4).Value = 2014
$ws.Cells.Item(14, 5).Value = -0.1259279434590921

$ws.Cells.Item(15, 1).Value = 41774
$ws.Cells.Item(15, 2).Value = 2014
$ws.Cells.Item(15, 3).Value = -0.1446844164011307
$ws.Cells.Item(15, 4).Value = 2015
$ws.Cells.Item(15, 5).Value = -0.03923323971219972

$ws.Cells.Item(16, 1).Value = 41957
$ws.Cells.Item(16, 2).Value = 2014
$ws.Cells.Item(16, 3).Value = 0.1729981757035093
$ws.Cells.Item(16, 4).Value = 2015
$ws.Cells.Item(16, 5).Value = 0.1749537368921361

$ws.Cells.Item(17, 1).Value = 42137
$ws.Cells.Item(17, 2).Value = 2015
$ws.Cells.Item(17, 3).Value = -0.1588690085687849
$ws.Cells.Item(17, 4).Value = 2016
$ws.Cells.Item(17, 5).Value = -0.4617525814883283

$ws.Cells.Item(18, 1).Value = 42321
$ws.Cells.Item(18, 2).Value = 2015
$ws.Cells.Item(18, 3).Value = 0.09752710595589686
$ws.Cells.Item(18, 4).Value = 2016
$ws.Cells.Item(18, 5).Value = -0.001769149545471915

$ws.Cells.Item(19, 1).Value = 42503
$ws.Cells.Item(19, 2).Value = 2016
$ws.Cells.Item(19, 3).Value = -0.5438176183081733
$ws.Cells.Item(19, 4).Value = 2017
$ws.Cells.Item(19, 5).Value = 0.01247916696665019

$ws.Cells.Item(20, 1).Value = 42689
$ws.Cells.Item(20, 2).Value = 2016
$ws.Cells.Item(20, 3).Value = -0.5280591151586633
$ws.Cells.Item(20, 4).Value = 2017
$ws.Cells.Item(20, 5).Value = -0.05116199209030947

$ws.Cells.Item(21, 1).Value = 42867
$ws.Cells.Item(21, 2).Value = 2017
$ws.Cells.Item(21, 3).Value = -0.006876704825709012
$ws.Cells.Item(21, 4).Value = 2018
$ws.Cells.Item(21, 5).Value = -0.05033169102144353

$ws.Cells.Item(22, 1).Value = 43053
$ws.Cells.Item(22, 2).Value = 2017
$ws.Cells.Item(22, 3).Value = 0.07201851318385799
$ws.Cells.Item(22, 4).Value = 2018
$ws.Cells.Item(22, 5).Value = 0.2735900898381383

$ws.Cells.Item(23, 1).Value = 43145
$ws.Cells.Item(23, 2).Value = 2018
$ws.Cells.Item(23, 3).Value = 0.6840863075407766
$ws.Cells.Item(23, 4).Value = 2019
$ws.Cells.Item(23, 5).Value = -0.0399940003999788

$ws.Cells.Item(24, 1).Value = 43235
$ws.Cells.Item(24, 2).Value = 2018
$ws.Cells.Item(24, 3).Value = 0.6718983809452572
$ws.Cells.Item(24, 4).Value = 2019
$ws.Cells.Item(24, 5).Value = 0.08750765859864007

$ws.Cells.Item(25, 1).Value = 43326
$ws.Cells.Item(25, 2).Value = 2018
$ws.Cells.Item(25, 3).Value = 0.4580297089606056
$ws.Cells.Item(25, 4).Value = 2019
$ws.Cells.Item(25, 5).Value = -0.09562311357415032

$ws.Cells.Item(26, 1).Value = 43418
$ws.Cells.Item(26, 2).Value = 2018
$ws.Cells.Item(26, 3).Value = 0.3727661260635617
$ws.Cells.Item(26, 4).Value = 2019
$ws.Cells.Item(26, 5).Value = -0.9505847809128332

$ws.Cells.Item(27, 1).Value = 43510
$ws.Cells.Item(27, 2).Value = 2019
$ws.Cells.Item(27, 3).Value = -0.8803054679952349
$ws.Cells.Item(27, 4).Value = 2020
$ws.Cells.Item(27, 5).Value = -0.2797061371759946

$ws.Cells.Item(28, 1).Value = 43600
$ws.Cells.Item(28, 2).Value = 2019
$ws.Cells.Item(28, 3).Value = -0.6347897325981511
$ws.Cells.Item(28, 4).Value = 2020
$ws.Cells.Item(28, 5).Value = 0.01241557525981651

$ws.Cells.Item(29, 1).Value = 43691
$ws.Cells.Item(29, 2).Value = 2019
$ws.Cells.Item(29, 3).Value = -0.8089889044073151
$ws.Cells.Item(29, 4).Value = 2020
$ws.Cells.Item(29, 5).Value = -0.2234428210501016

$ws.Cells.Item(30, 1).Value = 43783
$ws.Cells.Item(30, 2).Value = 2019
$ws.Cells.Item(30, 3).Value = -0.801759526476209
$ws.Cells.Item(30, 4).Value = 2020
$ws.Cells.Item(30, 5).Value = 0.047674034857903

$ws.Cells.Item(31, 1).Value = 43875
$ws.Cells.Item(31, 2).Value = 2020
$ws.Cells.Item(31, 3).Value = -0.3388987799285426
$ws.Cells.Item(31, 4).Value = 2021
$ws.Cells.Item(31, 5).Value = -0.1599040255974349

$ws.Cells.Item(32, 1).Value = 43966
$ws.Cells.Item(32, 2).Value = 2020
$ws.Cells.Item(32, 3).Value = -0.7158018152081724
$ws.Cells.Item(32, 4).Value = 2021
$ws.Cells.Item(32, 5).Value = -0.76158050880345

$ws.Cells.Item(33, 1).Value = 44068
$ws.Cells.Item(33, 2).Value = 2020
$ws.Cells.Item(33, 3).Value = -1.503583188367719
$ws.Cells.Item(33, 4).Value = 2021
$ws.Cells.Item(33, 5).Value = -0.8119518419852034

$ws.Cells.Item(34, 1).Value = 44159
$ws.Cells.Item(34, 2).Value = 2020
$ws.Cells.Item(34, 3).Value = -1.103489789942047
$ws.Cells.Item(34, 4).Value = 2021
$ws.Cells.Item(34, 5).Value = 1.605918384453009

$ws.Cells.Item(35, 1).Value = 44251
$ws.Cells.Item(35, 2).Value = 2021
$ws.Cells.Item(35, 3).Value = 2.893881462220338
$ws.Cells.Item(35, 4).Value = 2022
$ws.Cells.Item(35, 5).Value = 1.049469534781022

$ws.Cells.Item(36, 1).Value = 44341
$ws.Cells.Item(36, 2).Value = 2021
$ws.Cells.Item(36, 3).Value = 1.22331349480691
$ws.Cells.Item(36, 4).Value = 2022
$ws.Cells.Item(36, 5).Value = -0.5204428773059266

$ws.Cells.Item(37, 1).Value = 44432
$ws.Cells.Item(37, 2).Value = 2021
$ws.Cells.Item(37, 3).Value = 1.067534122491809
$ws.Cells.Item(37, 4).Value = 2022
$ws.Cells.Item(37, 5).Value = -0.2393417156003941

$ws.Cells.Item(38, 1).Value = 44525
$ws.Cells.Item(38, 2).Value = 2021
$ws.Cells.Item(38, 3).Value = 0.9704846793491928
$ws.Cells.Item(38, 4).Value = 2022
$ws.Cells.Item(38, 5).Value = -0.8255212498362474

$ws.Cells.Item(39, 1).Value = 44617
$ws.Cells.Item(39, 2).Value = 2022
$ws.Cells.Item(39, 3).Value = -0.4754331870696404
$ws.Cells.Item(39, 4).Value = 2023
$ws.Cells.Item(39, 5).Value = -0.3593817194810001

$ws.Cells.Item(40, 1).Value = 44706
$ws.Cells.Item(40, 2).Value = 2022
$ws.Cells.Item(40, 3).Value = -1.508346016334061
$ws.Cells.Item(40, 4).Value = 2023
$ws.Cells.Item(40, 5).Value = -1.034042971854776

$ws.Cells.Item(41, 1).Value = 44798
$ws.Cells.Item(41, 2).Value = 2022
$ws.Cells.Item(41, 3).Value = -0.9795431199870586
$ws.Cells.Item(41, 4).Value = 2023
$ws.Cells.Item(41, 5).Value = -0.476291821405983

$ws.Cells.Item(42, 1).Value = 44890
$ws.Cells.Item(42, 2).Value = 2022
$ws.Cells.Item(42, 3).Value = -0.7009264669202708
$ws.Cells.Item(42, 4).Value = 2023
$ws.Cells.Item(42, 5).Value = 0.6624163082313173

$ws.Cells.Item(43, 1).Value = 44981
$ws.Cells.Item(43, 2).Value = 2023
$ws.Cells.Item(43, 3).Value = 0.01219381132999686
$ws.Cells.Item(43, 4).Value = 2024
$ws.Cells.Item(43, 5).Value = -0.1052155539149613

$ws.Cells.Item(44, 1).Value = 45071
$ws.Cells.Item(44, 2).Value = 2023
$ws.Cells.Item(44, 3).Value = 0.2857158074419441
$ws.Cells.Item(44, 4).Value = 2024
$ws.Cells.Item(44, 5).Value = -0.2434890887128005

$ws.Cells.Item(45, 1).Value = 45163
$ws.Cells.Item(45, 2).Value = 2023
$ws.Cells.Item(45, 3).Value = 0.1829021030556488
$ws.Cells.Item(45, 4).Value = 2024
$ws.Cells.Item(45, 5).Value = -0.3407920741581916

$ws.Cells.Item(46, 1).Value = 45254
$ws.Cells.Item(46, 2).Value = 2023
$ws.Cells.Item(46, 3).Value = 0.3928252664241905
$ws.Cells.Item(46, 4).Value = 2024
$ws.Cells.Item(46, 5).Value = 0.302295480375836

$ws.Cells.Item(47, 1).Value = 45345
$ws.Cells.Item(47, 2).Value = 2024
$ws.Cells.Item(47, 3).Value = 0.2614828632771848
$ws.Cells.Item(47, 4).Value = 2025
$ws.Cells.Item(47, 5).Value = 0.1003260477960621

$ws.Cells.Item(48, 1).Value = 45436
$ws.Cells.Item(48, 2).Value = 2024
$ws.Cells.Item(48, 3).Value = 1.050311853611596
$ws.Cells.Item(48, 4).Value = 2025
$ws.Cells.Item(48, 5).Value = 0.6547560647617745

$ws.Cells.Item(49, 1).Value = 45534
$ws.Cells.Item(49, 2).Value = 2024
$ws.Cells.Item(49, 3).Value = 0.6979546684258597
$ws.Cells.Item(49, 4).Value = 2025
$ws.Cells.Item(49, 5).Value = 0.1438580449789839

$ws.Cells.Item(50, 1).Value = 45618
$ws.Cells.Item(50, 2).Value = 2024
$ws.Cells.Item(50, 3).Value = 0.3224026462283813
$ws.Cells.Item(50, 4).Value = 2025
$ws.Cells.Item(50, 5).Value = -0.7618983399156787

$ws.Cells.Item(51, 1).Value = 45713
$ws.Cells.Item(51, 2).Value = 2025
$ws.Cells.Item(51, 3).Value = -2.451829860015453
$ws.Cells.Item(51, 4).Value = 2026
$ws.Cells.Item(51, 5).Value = -0.5076775133253331

$ws.Cells.Item(52, 1).Value = 45800
$ws.Cells.Item(52, 2).Value = 2025
$ws.Cells.Item(52, 3).Value = -1.467296258526263
$ws.Cells.Item(52, 4).Value = 2026
$ws.Cells.Item(52, 5).Value = -0.02957675682233596

$ws.Cells.Item(53, 1).Value = 45891
$ws.Cells.Item(53, 2).Value = 2025
$ws.Cells.Item(53, 3).Value = -2.11737366557071
$ws.Cells.Item(53, 4).Value = 2026
$ws.Cells.Item(53, 5).Value = -0.9529168788645181

